# Daily update: append the latest daily_profile_metrics row (row 44) and
# correct the stored extraction_datetime precision on the previous row (43).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Tiny precision fix on the existing last row (row 43) ---
# Only the G (extraction_datetime serial) value changes; nothing else about
# this row is touched so its existing number formatting is left alone.
$ws.Range("G43").Value = 45789.42137028935

# --- Append the new day's scrape as row 44 ---
$bio = "ʜᴇʟᴘɪɴɢ ʏᴏᴜ ᴍᴏᴠᴇ ғʀᴏᴍ ʜᴇsɪᴛᴀᴛɪᴏɴ ᴛᴏ ᴄʀᴇᴀᴛɪᴏɴ`n✨ | ✧ 𝗹𝗼𝗰𝘀 ✧ (𝘀𝗲𝗹𝗳) 𝗹𝗶𝗯𝗲𝗿𝗮𝘁𝗶𝗼𝗻 ✧ 𝗹𝗶𝗳𝗲𝘀𝘁𝘆𝗹𝗲 ✧`n🪴 | 71 ʟᴏᴄs est. on 07.20.23`n📍 | ʜᴏᴜsᴛᴏɴ, ᴛx"

$ws.Range("A44").Value = $bio
$ws.Range("B44").Value = 2984
$ws.Range("C44").Value = 238
$ws.Range("D44").Value = 174
$ws.Range("E44").Value = "https://scontent-hou1-1.xx.fbcdn.net/v/t51.2885-15/481266977_997353345602937_1719041919639027270_n.jpg?_nc_cat=106&ccb=1-7&_nc_sid=7d201b&_nc_ohc=98Z30Y-Fn84Q7kNvwFG_Vub&_nc_oc=AdlBhQs5yoAx2HX9ehwQEDxKusa4_dxclDLJsZm2zFF7LvHd2ATdglLLxKuRHjvPxDybk4J6fG9sXmCu5a9P_MPF&_nc_zt=23&_nc_ht=scontent-hou1-1.xx&edm=AL-3X8kEAAAA&oh=00_AfJ-Y6F611Y1W1BQhsDBLCP4JvvrgYAWu8LlVQrb2jg1fw&oe=68295BAE"
$ws.Range("F44").Value = "17841461458191255"
$ws.Range("G44").Value = 45790.59407874689
$ws.Range("G44").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H44").Value = "2025-05-13"
$ws.Range("I44").Value = 2025
$ws.Range("J44").Value = "May"
$ws.Range("K44").Value = 13
$ws.Range("L44").Value = "14:15:28"
